$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late") to make room for a
# "Variable Instalments" column — shifts Late/heading/Outstanding one
# column to the right (N->O, O->P, P->Q).
$ws.Columns.Item(14).Insert()

# Match the width Excel would carry over from the neighbouring column (M)
# when inserting a new column.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth()

# Make "Repayment schedule" the active sheet/tab and select cell R6 on it,
# matching the saved view state of the edited workbook.
$ws.Activate()
$ws.Range("R6").Select()

Write-Host "done"
